$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AC, AD, AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy style from an existing header cell (e.g. AB1) to the new header cells
$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in team record data for rows 2-50
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 29).Value = 84   # AC
    $ws.Cells.Item($r, 30).Value = 78   # AD
    $ws.Cells.Item($r, 31).Value = 0    # AE
}
